$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B12").Value = 3.0
$ws.Range("B20").Value = 3.0
$ws.Range("B37").Value = 3.0
$ws.Range("B39").Value = 5.0
$ws.Range("B42").Value = 1.0
$ws.Range("B43").Value = 1.0
$ws.Range("B46").Value = 1.0
$ws.Range("B47").Value = 3.0
$ws.Range("B48").Value = 1.0
$ws.Range("B54").Value = 4.0
$ws.Range("B55").Value = 2.0
$ws.Range("B58").Value = 3.0
$ws.Range("B61").Value = 3.0
$ws.Range("B70").Value = 3.0
$ws.Range("B71").Value = 3.0
$ws.Range("B79").Value = 3.0
$ws.Range("B81").Value = 2.0
$ws.Range("B82").Value = 1.0
$ws.Range("B90").Value = 3.0
$ws.Range("B92").Value = 4.0
$ws.Range("B94").Value = 3.0
$ws.Range("B97").Value = 1.0
$ws.Range("B101").Value = 1.0
